$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.570.19'
$ws.Range("E2").Value = '  +0.11%  '
# Row 3
$ws.Range("D3").Value = '1.754.61'
$ws.Range("E3").Value = '  +0.15%  '
# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.33%  '
# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '324.29'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4513'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.21%  '
# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3572'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.64%  '
# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07471'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.39%  '
# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '41.46'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.82%  '
# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.087'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.26%  '
# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.27%  '
# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '20.79'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.53%  '
# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.991'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.70%  '
# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.165'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.45%  '
# Row 16
$ws.Range("D16").Value = '1.754.08'
$ws.Range("E16").Value = '  -0.28%  '
# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '94.52'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.92%  '
# Row 18
$ws.Range("E18").Value = '  -0.62%  '
# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06388'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.37%  '
# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.17%  '
# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '17.11'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.38%  '
# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.747'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.78%  '
# Row 23
$ws.Range("D23").Value = '27.616.57'
$ws.Range("E23").Value = '  +0.10%  '
# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '11.20'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '
# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.082'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.81%  '
# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '165.70'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.92%  '
# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '20.20'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.94%  '
# Row 28
$ws.Range("D28").Value = '1.955.66'
$ws.Range("E28").Value = '  -0.13%  '
# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.134'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.45%  '
# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '125.64'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.17%  '
# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.085'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.38%  '
# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.09172'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.45%  '
# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.649'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.26%  '
# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.510'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.46%  '
# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.02287'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.48%  '
# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '11.72'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.62%  '
# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.2092'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.49%  '
# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.06003'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.50%  '
# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.6289'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '
# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '4.926'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '
# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.182'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.91%  '
# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.390'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.10%  '
# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '7.792'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.12%  '
# Row 44
$ws.Range("E44").Value = '  -0.60%  '
# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.715'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.09%  '
# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.5862'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.14%  '
# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '121.85'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.25%  '
# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.935'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.27%  '
# Row 49
$ws.Range("E49").Value = '  +0.23%  '
# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.131'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.04%  '
# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '71.50'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.79%  '
